$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly row at row 6 (pushes old rows 6-19 down to 7-20)
$ws.Rows.Item(6).Insert()

# Insert another brand-new weekly row at row 10 (in the new numbering),
# pushing the rows that are now 10-20 further down to 11-21
$ws.Rows.Item(10).Insert()

# --- Populate new row 6 ---
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(6, 3).Value = "Coquimbo"
$ws.Cells.Item(6, 4).Value = 44630
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100101
$ws.Cells.Item(6, 8).Value = "Berries"
$ws.Cells.Item(6, 9).Value = 100101001
$ws.Cells.Item(6, 10).Value = "Arándano (blue)"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 240
$ws.Cells.Item(6, 14).Value = 6000
$ws.Cells.Item(6, 15).Value = 6500
$ws.Cells.Item(6, 16).Value = 6250
$ws.Cells.Item(6, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(6, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(6, 19).Value = 3125
$ws.Cells.Item(6, 20).Value = 2

# --- Populate new row 10 ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44435
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100101
$ws.Cells.Item(10, 8).Value = "Berries"
$ws.Cells.Item(10, 9).Value = 100101001
$ws.Cells.Item(10, 10).Value = "Arándano (blue)"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 400
$ws.Cells.Item(10, 14).Value = 19500
$ws.Cells.Item(10, 15).Value = 20000
$ws.Cells.Item(10, 16).Value = 19750
$ws.Cells.Item(10, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 19).Value = 9875
$ws.Cells.Item(10, 20).Value = 2
